$d = $word.ActiveDocument

# Legacy VML shape-id renumbering (o:spid/_x0000_iNNNN) applied when the
# ActiveX/OLE "team report" controls were re-inserted. Each w:object run
# carries a v:shape id="_x0000_iNNNN" plus a matching w:control
# w:shapeid="_x0000_iNNNN" that must be updated together.
$idMap = @{
    "_x0000_i1095" = "_x0000_i1026";
    "_x0000_i1096" = "_x0000_i1028";
    "_x0000_i1097" = "_x0000_i1030";
    "_x0000_i1098" = "_x0000_i1032";
    "_x0000_i1099" = "_x0000_i1034";
    "_x0000_i1100" = "_x0000_i1036";
    "_x0000_i1101" = "_x0000_i1038";
    "_x0000_i1102" = "_x0000_i1040";
    "_x0000_i1103" = "_x0000_i1042";
    "_x0000_i1104" = "_x0000_i1044";
    "_x0000_i1105" = "_x0000_i1046";
    "_x0000_i1106" = "_x0000_i1048";
    "_x0000_i1107" = "_x0000_i1050";
    "_x0000_i1093" = "_x0000_i1052";
}

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $xml = $p.Range.WordOpenXML

    foreach ($oldId in $idMap.Keys) {
        if ($xml -like "*$oldId*") {
            $newId = $idMap[$oldId]
            $newXml = $xml.Replace($oldId, $newId)

            # Range.WordOpenXML synthesizes w14:paraId/w14:textId for any
            # paragraph that doesn't already carry one (none do, in this
            # document) -- strip that synthetic stamp back out before
            # InsertXML so re-inserting the fragment doesn't newly persist
            # attributes the original paragraph never had.
            $newXml = $newXml -replace ' w14:paraId="[0-9A-Fa-f]{8}" w14:textId="[0-9A-Fa-f]{8}"', ''

            $p.Range.InsertXML($newXml)
            break
        }
    }
}
